# Update each slide's auto-advance ("Advance Slide After") timing.
#
# PowerPoint exposes this as SlideShowTransition.AdvanceTime, a value in
# seconds; on disk it is stored (in milliseconds) as the advTm attribute
# of the slide's <p:transition> element.

$p = $ppt.ActivePresentation

$slideIndexes  = @(1, 2, 3, 4, 5, 6, 7, 8)
$advanceTimes  = @(23.798, 64.179, 6.205, 86.211, 31.281, 13.541, 21.424, 14.162)

for ($i = 0; $i -lt $slideIndexes.Count; $i++) {
    $slide = $p.Slides.Item($slideIndexes[$i])
    $slide.SlideShowTransition.AdvanceTime = $advanceTimes[$i]
}
